# Map_to_Geojson_Gantt.xlsx - v0.0.6 "Improved all files and added other tests"
# Translates the Gantt task list from Italian to English, bumps the version
# label, tweaks a couple of progress percentages, and updates the view state
# (zoom level + selection) to match the author's final save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / version banner ---------------------------------------------
$ws.Range("C2").Value = "Map to .geojson converter"
$ws.Range("C4").Value = "Gantt v3.0.0"

# --- Task names (column C) translated to English -------------------------
$ws.Range("C7").Value  = "Research & Prototyping"
$ws.Range("C8").Value  = "Feasibility Study"
$ws.Range("C9").Value  = "Gantt Prototype"
$ws.Range("C10").Value = "SVG & GeoJSON study"
$ws.Range("C11").Value = "Dataset Study"
$ws.Range("C12").Value = "SVG to GeoJSON protype"
$ws.Range("C13").Value = "Computer Vision Development"
$ws.Range("C14").Value = "Dataset Research"
$ws.Range("C15").Value = "Libraries Research"
$ws.Range("C16").Value = "Initial Implementation"
$ws.Range("C17").Value = "More Features"
$ws.Range("C18").Value = "UI Development"
$ws.Range("C19").Value = "Figma design"
$ws.Range("C20").Value = "Front-end implementation"
$ws.Range("C21").Value = "Testing & Release"
$ws.Range("C22").Value = "Testing"
$ws.Range("C23").Value = "Testing + UI"
$ws.Range("C24").Value = "Beta release and user testing"
$ws.Range("C26").Value = "Brand study"
$ws.Range("C27").Value = "Market entry planning"
$ws.Range("C28").Value = "Marketing planning"

# --- Collaborators label ---------------------------------------------------
$ws.Range("C32").Value = "Collaborators"

# --- Progress percentage tweaks ------------------------------------------
$ws.Range("E10").Value = 0.4
$ws.Range("E26").Value = 0.15

# --- View state: zoom + selection -----------------------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("C28").Select()
